$wb = $excel.ActiveWorkbook

# Rename the first sheet per the diff
$wb.Worksheets.Item("Emmer and Einkorn").Name = "Wild Emmer and Einkorn"

$wsWildEmmer = $wb.Worksheets.Item("Wild Emmer and Einkorn")
$wsWildEmmer.Range("B2").Value = 0.3721544792316152
$wsWildEmmer.Range("C2").Value = 0.1043231789536628
$wsWildEmmer.Range("B3").Value = 0.3508204786128634
$wsWildEmmer.Range("C3").Value = -0.7455931581521854
$wsWildEmmer.Range("B4").Value = 0.3526999915615989
$wsWildEmmer.Range("C4").Value = 0.6552583507526101
$wsWildEmmer.Range("B5").Value = 0.4583657748124931
$wsWildEmmer.Range("C5").Value = -0.003291137479821626
$wsWildEmmer.Range("B6").Value = 0.4399877922076786
$wsWildEmmer.Range("C6").Value = -0.05146209735416027
$wsWildEmmer.Range("B7").Value = 0.4586287284011868
$wsWildEmmer.Range("C7").Value = 0.0344213824333387

$wsEinkorn = $wb.Worksheets.Item("Einkorn")
$wsEinkorn.Range("B2").Value = -0.4087576322602193
$wsEinkorn.Range("C2").Value = 0.1059363649697134
$wsEinkorn.Range("B3").Value = -0.2383853422796498
$wsEinkorn.Range("C3").Value = -0.8704158396480776
$wsEinkorn.Range("B4").Value = -0.3798801081040958
$wsEinkorn.Range("C4").Value = 0.4719196994013413
$wsEinkorn.Range("B5").Value = -0.4600816232039395
$wsEinkorn.Range("B6").Value = -0.4537686903103582
$wsEinkorn.Range("C6").Value = -0.05827149681265453
$wsEinkorn.Range("B7").Value = -0.4628170327527977
$wsEinkorn.Range("C7").Value = 0.06091995927444105

$wsEmmer = $wb.Worksheets.Item("Emmer")
$wsEmmer.Range("B2").Value = -0.4163999084396875
$wsEmmer.Range("C2").Value = -0.3013862094196741
$wsEmmer.Range("B3").Value = -0.3087437282723052
$wsEmmer.Range("C3").Value = 0.7228606580623405
$wsEmmer.Range("B4").Value = -0.3866476059230858
$wsEmmer.Range("C4").Value = -0.4945255840539674
$wsEmmer.Range("B5").Value = -0.4577488248650109
$wsEmmer.Range("C5").Value = -0.05135376020091403
$wsEmmer.Range("B6").Value = -0.4031956695964494
$wsEmmer.Range("C6").Value = 0.3672247980680073
$wsEmmer.Range("B7").Value = -0.4579206489024669
$wsEmmer.Range("C7").Value = -0.06776435280102949

$wsBarley = $wb.Worksheets.Item("Barley")
$wsBarley.Range("B2").Value = 0.3903008835525389
$wsBarley.Range("C2").Value = -0.4128639341577707
$wsBarley.Range("B3").Value = 0.2040553963386507
$wsBarley.Range("C3").Value = 0.731065325432493
$wsBarley.Range("B4").Value = 0.3595474228337249
$wsBarley.Range("B5").Value = 0.4816467844222868
$wsBarley.Range("C5").Value = 0.09964651911847838
$wsBarley.Range("C6").Value = 0.1880429904741636
$wsBarley.Range("B7").Value = 0.4806938224648517
$wsBarley.Range("C7").Value = 0.1090358961104006

